# Auto-generated: update stat values in row 4-13 per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 0.3
$ws.Range("E4").Value = 0.18
$ws.Range("F4").Value = 0.028
$ws.Range("G4").Value = 0.167
$ws.Range("H4").Value = 0.214
$ws.Range("J4").Value = 0.157
$ws.Range("K4").Value = 0.329
$ws.Range("L4").Value = 0.099
$ws.Range("M4").Value = 0.314
$ws.Range("N4").Value = 0.266
$ws.Range("O4").Value = 0.02
$ws.Range("P4").Value = 0.143
$ws.Range("Q4").Value = 0.512
$ws.Range("R4").Value = 0.22
$ws.Range("S4").Value = 0.469
$ws.Range("T4").Value = 0.262
$ws.Range("U4").Value = 0.085
$ws.Range("W4").Value = 0.251
$ws.Range("X4").Value = 0.043
$ws.Range("Y4").Value = 0.206
$ws.Range("Z4").Value = 0.462
$ws.Range("AA4").Value = 0.125
$ws.Range("AB4").Value = 0.354
$ws.Range("AC4").Value = 0.121
$ws.Range("AD4").Value = 0.006
$ws.Range("AE4").Value = 0.08
$ws.Range("AF4").Value = 0.71
$ws.Range("AG4").Value = 0.104
$ws.Range("AH4").Value = 0.323
$ws.Range("AI4").Value = 0.667
$ws.Range("AJ4").Value = 0.163
$ws.Range("AK4").Value = 0.403
$ws.Range("AL4").Value = 0.687
$ws.Range("AM4").Value = 0.113
$ws.Range("AN4").Value = 0.336
$ws.Range("AO4").Value = 0.688
# Row 5
$ws.Range("B5").Value = 0.833
$ws.Range("C5").Value = 0.139
$ws.Range("D5").Value = 0.373
$ws.Range("E5").Value = 0.714
$ws.Range("F5").Value = 0.204
$ws.Range("G5").Value = 0.452
$ws.Range("H5").Value = 0.857
$ws.Range("I5").Value = 0.122
$ws.Range("J5").Value = 0.35
$ws.Range("K5").Value = 0.643
$ws.Range("L5").Value = 0.23
$ws.Range("M5").Value = 0.479
$ws.Range("N5").Value = 0.833
$ws.Range("O5").Value = 0.139
$ws.Range("P5").Value = 0.373
$ws.Range("Q5").Value = 0.571
$ws.Range("R5").Value = 0.245
$ws.Range("S5").Value = 0.495
$ws.Range("T5").Value = 0.548
$ws.Range("U5").Value = 0.248
$ws.Range("V5").Value = 0.498
$ws.Range("W5").Value = 0.762
$ws.Range("X5").Value = 0.181
$ws.Range("Y5").Value = 0.426
$ws.Range("Z5").Value = 0.857
$ws.Range("AA5").Value = 0.122
$ws.Range("AB5").Value = 0.35
$ws.Range("AC5").Value = 0.738
$ws.Range("AD5").Value = 0.193
$ws.Range("AE5").Value = 0.44
$ws.Range("AF5").Value = 0.952
$ws.Range("AG5").Value = 0.045
$ws.Range("AH5").Value = 0.213
$ws.Range("AI5").Value = 0.786
$ws.Range("AJ5").Value = 0.168
$ws.Range("AK5").Value = 0.41
$ws.Range("AL5").Value = 0.929
$ws.Range("AM5").Value = 0.066
$ws.Range("AN5").Value = 0.258
$ws.Range("AO5").Value = 0.889
# Row 6
$ws.Range("B6").Value = 0.441
$ws.Range("E6").Value = 0.288
$ws.Range("H6").Value = 0.342
$ws.Range("K6").Value = 0.435
$ws.Range("N6").Value = 0.403
$ws.Range("Q6").Value = 0.54
$ws.Range("T6").Value = 0.355
$ws.Range("W6").Value = 0.378
$ws.Range("Z6").Value = 0.6
$ws.Range("AC6").Value = 0.208
$ws.Range("AF6").Value = 0.813
$ws.Range("AI6").Value = 0.722
$ws.Range("AL6").Value = 0.79
$ws.Range("AO6").Value = 0.775
# Row 7
$ws.Range("B7").Value = 0.615
$ws.Range("E7").Value = 0.448
$ws.Range("H7").Value = 0.535
$ws.Range("K7").Value = 0.54
$ws.Range("N7").Value = 0.584
$ws.Range("Q7").Value = 0.558
$ws.Range("T7").Value = 0.45
$ws.Range("W7").Value = 0.542
$ws.Range("Z7").Value = 0.732
$ws.Range("AC7").Value = 0.365
$ws.Range("AF7").Value = 0.891
$ws.Range("AI7").Value = 0.759
$ws.Range("AL7").Value = 0.868
$ws.Range("AO7").Value = 0.839
# Row 8
$ws.Range("B8").Value = 0.759
$ws.Range("C8").Value = 0.143
$ws.Range("D8").Value = 0.378
$ws.Range("E8").Value = 0.603
$ws.Range("G8").Value = 0.426
$ws.Range("H8").Value = 0.747
$ws.Range("I8").Value = 0.132
$ws.Range("J8").Value = 0.364
$ws.Range("K8").Value = 0.569
$ws.Range("N8").Value = 0.752
$ws.Range("O8").Value = 0.141
$ws.Range("P8").Value = 0.376
$ws.Range("Q8").Value = 0.545
$ws.Range("R8").Value = 0.231
$ws.Range("S8").Value = 0.481
$ws.Range("T8").Value = 0.477
$ws.Range("V8").Value = 0.458
$ws.Range("W8").Value = 0.685
$ws.Range("X8").Value = 0.17
$ws.Range("Y8").Value = 0.413
$ws.Range("Z8").Value = 0.789
$ws.Range("AA8").Value = 0.127
$ws.Range("AB8").Value = 0.356
$ws.Range("AC8").Value = 0.615
$ws.Range("AE8").Value = 0.426
$ws.Range("AF8").Value = 0.879
$ws.Range("AG8").Value = 0.062
$ws.Range("AH8").Value = 0.248
$ws.Range("AI8").Value = 0.777
$ws.Range("AJ8").Value = 0.168
$ws.Range("AK8").Value = 0.41
$ws.Range("AL8").Value = 0.893
$ws.Range("AM8").Value = 0.073
$ws.Range("AN8").Value = 0.27
$ws.Range("AO8").Value = 0.85
# Row 9
$ws.Range("B9").Value = 0.667
$ws.Range("C9").Value = 0.222
$ws.Range("D9").Value = 0.471
$ws.Range("E9").Value = 0.476
$ws.Range("H9").Value = 0.619
$ws.Range("I9").Value = 0.236
$ws.Range("J9").Value = 0.486
$ws.Range("K9").Value = 0.476
$ws.Range("N9").Value = 0.643
$ws.Range("O9").Value = 0.23
$ws.Range("P9").Value = 0.479
$ws.Range("Q9").Value = 0.5
$ws.Range("T9").Value = 0.381
$ws.Range("U9").Value = 0.236
$ws.Range("V9").Value = 0.486
$ws.Range("W9").Value = 0.571
$ws.Range("X9").Value = 0.245
$ws.Range("Y9").Value = 0.495
$ws.Range("Z9").Value = 0.69
$ws.Range("AA9").Value = 0.214
$ws.Range("AB9").Value = 0.462
$ws.Range("AC9").Value = 0.5
$ws.Range("AF9").Value = 0.762
$ws.Range("AG9").Value = 0.181
$ws.Range("AH9").Value = 0.426
$ws.Range("AI9").Value = 0.762
$ws.Range("AJ9").Value = 0.181
$ws.Range("AK9").Value = 0.426
$ws.Range("AL9").Value = 0.833
$ws.Range("AM9").Value = 0.139
$ws.Range("AN9").Value = 0.373
$ws.Range("AO9").Value = 0.786
# Row 10
$ws.Range("B10").Value = 0.786
$ws.Range("C10").Value = 0.168
$ws.Range("D10").Value = 0.41
$ws.Range("E10").Value = 0.643
$ws.Range("F10").Value = 0.23
$ws.Range("G10").Value = 0.479
$ws.Range("H10").Value = 0.786
$ws.Range("I10").Value = 0.168
$ws.Range("J10").Value = 0.41
$ws.Range("K10").Value = 0.643
$ws.Range("L10").Value = 0.23
$ws.Range("M10").Value = 0.479
$ws.Range("N10").Value = 0.81
$ws.Range("O10").Value = 0.154
$ws.Range("P10").Value = 0.393
$ws.Range("Q10").Value = 0.571
$ws.Range("R10").Value = 0.245
$ws.Range("S10").Value = 0.495
$ws.Range("T10").Value = 0.548
$ws.Range("U10").Value = 0.248
$ws.Range("V10").Value = 0.498
$ws.Range("W10").Value = 0.762
$ws.Range("X10").Value = 0.181
$ws.Range("Y10").Value = 0.426
$ws.Range("Z10").Value = 0.857
$ws.Range("AA10").Value = 0.122
$ws.Range("AB10").Value = 0.35
$ws.Range("AC10").Value = 0.619
$ws.Range("AD10").Value = 0.236
$ws.Range("AE10").Value = 0.486
$ws.Range("AF10").Value = 0.952
$ws.Range("AG10").Value = 0.045
$ws.Range("AH10").Value = 0.213
$ws.Range("AI10").Value = 0.786
$ws.Range("AJ10").Value = 0.168
$ws.Range("AK10").Value = 0.41
$ws.Range("AL10").Value = 0.929
$ws.Range("AM10").Value = 0.066
$ws.Range("AN10").Value = 0.258
$ws.Range("AO10").Value = 0.889
# Row 11
$ws.Range("B11").Value = 0.833
$ws.Range("C11").Value = 0.139
$ws.Range("D11").Value = 0.373
$ws.Range("E11").Value = 0.714
$ws.Range("F11").Value = 0.204
$ws.Range("G11").Value = 0.452
$ws.Range("H11").Value = 0.857
$ws.Range("I11").Value = 0.122
$ws.Range("J11").Value = 0.35
$ws.Range("K11").Value = 0.643
$ws.Range("L11").Value = 0.23
$ws.Range("M11").Value = 0.479
$ws.Range("N11").Value = 0.833
$ws.Range("O11").Value = 0.139
$ws.Range("P11").Value = 0.373
$ws.Range("Q11").Value = 0.571
$ws.Range("R11").Value = 0.245
$ws.Range("S11").Value = 0.495
$ws.Range("T11").Value = 0.548
$ws.Range("U11").Value = 0.248
$ws.Range("V11").Value = 0.498
$ws.Range("W11").Value = 0.762
$ws.Range("X11").Value = 0.181
$ws.Range("Y11").Value = 0.426
$ws.Range("Z11").Value = 0.857
$ws.Range("AA11").Value = 0.122
$ws.Range("AB11").Value = 0.35
$ws.Range("AC11").Value = 0.667
$ws.Range("AD11").Value = 0.222
$ws.Range("AE11").Value = 0.471
$ws.Range("AF11").Value = 0.952
$ws.Range("AG11").Value = 0.045
$ws.Range("AH11").Value = 0.213
$ws.Range("AI11").Value = 0.786
$ws.Range("AJ11").Value = 0.168
$ws.Range("AK11").Value = 0.41
$ws.Range("AL11").Value = 0.929
$ws.Range("AM11").Value = 0.066
$ws.Range("AN11").Value = 0.258
$ws.Range("AO11").Value = 0.889
# Row 12
$ws.Range("B12").Value = 1.343
$ws.Range("C12").Value = 0.625
$ws.Range("D12").Value = 0.791
$ws.Range("E12").Value = 1.633
$ws.Range("F12").Value = 1.032
$ws.Range("G12").Value = 1.016
$ws.Range("H12").Value = 1.556
$ws.Range("I12").Value = 1.191
$ws.Range("J12").Value = 1.091
$ws.Range("K12").Value = 1.407
$ws.Range("L12").Value = 0.538
$ws.Range("M12").Value = 0.733
$ws.Range("N12").Value = 1.343
$ws.Range("O12").Value = 0.511
$ws.Range("P12").Value = 0.715
$ws.Range("Z12").Value = 1.25
$ws.Range("AA12").Value = 0.299
$ws.Range("AB12").Value = 0.546
$ws.Range("AC12").Value = 2.032
$ws.Range("AD12").Value = 3.902
$ws.Range("AE12").Value = 1.975
$ws.Range("AF12").Value = 1.225
$ws.Range("AG12").Value = 0.224
$ws.Range("AH12").Value = 0.474
$ws.Range("AI12").Value = 1.03
$ws.Range("AJ12").Value = 0.029
$ws.Range("AK12").Value = 0.171
$ws.Range("AL12").Value = 1.103
$ws.Range("AM12").Value = 0.092
$ws.Range("AN12").Value = 0.303
$ws.Range("AO12").Value = 1.119
# Row 13
$ws.Range("B13").Value = 3.429
$ws.Range("C13").Value = 1.34
$ws.Range("D13").Value = 1.158
$ws.Range("E13").Value = 4.541
$ws.Range("F13").Value = 0.735
$ws.Range("G13").Value = 0.857
$ws.Range("H13").Value = 4.5
$ws.Range("I13").Value = 0.95
$ws.Range("J13").Value = 0.975
$ws.Range("K13").Value = 2.333
$ws.Range("L13").Value = 0.581
$ws.Range("M13").Value = 0.762
$ws.Range("N13").Value = 3.286
$ws.Range("O13").Value = 0.776
$ws.Range("P13").Value = 0.881
$ws.Range("Z13").Value = 2.833
$ws.Range("AA13").Value = 3.901
$ws.Range("AB13").Value = 1.975
$ws.Range("AC13").Value = 6.268
$ws.Range("AD13").Value = 3.026
$ws.Range("AE13").Value = 1.739
$ws.Range("AF13").Value = 1.643
$ws.Range("AG13").Value = 0.706
$ws.Range("AH13").Value = 0.84
$ws.Range("AI13").Value = 1.238
$ws.Range("AJ13").Value = 0.181
$ws.Range("AK13").Value = 0.426
$ws.Range("AL13").Value = 1.667
$ws.Range("AM13").Value = 0.794
$ws.Range("AN13").Value = 0.891
$ws.Range("AO13").Value = 1.516

Write-Output "updated 313 cells"
